$wb = $excel.ActiveWorkbook

# --- Sheet 1: "ATS Accuracy" ---
$ws1 = $wb.Worksheets.Item("ATS Accuracy")

# Row 2
$ws1.Range("B2").Value = 6
$ws1.Range("C2").Value = 88
$ws1.Range("D2").Value = 94
$ws1.Range("E2").Value = 93.59999999999999

# Row 3
$ws1.Range("B3").Value = 10
$ws1.Range("C3").Value = 68
$ws1.Range("D3").Value = 78
$ws1.Range("E3").Value = 87.2

# Row 4
$ws1.Range("C4").Value = 15
$ws1.Range("D4").Value = 19
$ws1.Range("E4").Value = 78.90000000000001

# Row 5
$ws1.Range("B5").Value = 2
$ws1.Range("C5").Value = 7
$ws1.Range("D5").Value = 9
$ws1.Range("E5").Value = 77.8

# --- Sheet 2: "Total Accuracy" ---
$ws2 = $wb.Worksheets.Item("Total Accuracy")

# Row 2
$ws2.Range("B2").Value = 5
$ws2.Range("C2").Value = 77
$ws2.Range("D2").Value = 82
$ws2.Range("E2").Value = 93.90000000000001

# Row 3
$ws2.Range("B3").Value = 8
$ws2.Range("C3").Value = 74
$ws2.Range("D3").Value = 82
$ws2.Range("E3").Value = 90.2

# Row 4
$ws2.Range("B4").Value = 3
$ws2.Range("C4").Value = 21
$ws2.Range("D4").Value = 24
$ws2.Range("E4").Value = 87.5

# Row 5
$ws2.Range("B5").Value = 6
$ws2.Range("C5").Value = 9
$ws2.Range("D5").Value = 15
$ws2.Range("E5").Value = 60

# Row 6
$ws2.Range("B6").Value = 3
$ws2.Range("C6").Value = 2
$ws2.Range("D6").Value = 5
$ws2.Range("E6").Value = 40
